$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 46060 to 46061 for rows 2-14
$ws.Range("C2:C14").Value = 46061
